$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New version row 7 gained a "Label on gateway" value of "?" (undefined label on gateway)
$cell = $ws.Range("D7")
$cell.Value = "?"

# Style: red font, centered horizontally
$cell.Font.Color = 255
$cell.HorizontalAlignment = -4108  # xlCenter

# Move the active selection to D7 (matches the saved selection in the sheet view)
$ws.Range("D7").Select()
